$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record as row 62, pushing the existing rows 62-96
# (and everything below) down by one row.
$ws.Range("A62:R62").EntireRow.Insert()

# Populate the newly inserted row with this week's values.
$ws.Range("A62").Value = 8
$ws.Range("B62").Value = "Terminal La Palmera de La Serena"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 44518
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = 100112040
$ws.Range("G62").Value = "Cilantro"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 3200
$ws.Range("K62").Value = 1300
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = 1400
$ws.Range("N62").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O62").Value = "Provincia del Elquí"
$ws.Range("P62").Value = 933
$ws.Range("Q62").Value = 1.5
$ws.Range("R62").Value = "Hortaliza"
